$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

# Header cell F1, styled like the other header cells (B1:E1)
$ws.Range("F1").Value = "time_taken"
$ws.Range("F1").Style = $ws.Range("E1").Style

# Data cells F2:F5 with timestamp values matching plain (unstyled) data cells
$ws.Range("F2").Value = "2021-10-05 13:38:57.776036"
$ws.Range("F3").Value = "2021-10-05 13:38:57.776050"
$ws.Range("F4").Value = "2021-10-05 13:38:57.776054"
$ws.Range("F5").Value = "2021-10-05 13:38:57.776058"
